$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 540, pushing existing rows 540-601 down to 541-602
$ws.Rows.Item(540).Insert()

# Populate the newly inserted row 540 with the new weekly price entry
$ws.Cells.Item(540, 1).Value2 = 4
$ws.Cells.Item(540, 2).Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(540, 3).Value2 = "Los Lagos"
$ws.Cells.Item(540, 4).Value2 = 45194
$ws.Cells.Item(540, 5).Value2 = 10
$ws.Cells.Item(540, 6).Value2 = 100112023
$ws.Cells.Item(540, 7).Value2 = "Brócoli"
$ws.Cells.Item(540, 8).Value2 = "Sin especificar"
$ws.Cells.Item(540, 9).Value2 = "Primera"
$ws.Cells.Item(540, 10).Value2 = 500
$ws.Cells.Item(540, 11).Value2 = 1400
$ws.Cells.Item(540, 12).Value2 = 1400
$ws.Cells.Item(540, 13).Value2 = 1400
$ws.Cells.Item(540, 14).Value2 = "$/unidad"
$ws.Cells.Item(540, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(540, 16).Value2 = 1400
$ws.Cells.Item(540, 17).Value2 = 1
$ws.Cells.Item(540, 18).Value2 = "Hortaliza"

# Keep the same date number format used by the rest of column D
$ws.Cells.Item(540, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
